# Horarios actualizados Línea 141 - 1266
# Updates "Última actualización" timestamp across all sheets, refreshes the
# LP1912 schedule rows (sheet 1), and appends a new arrival row.

$wb = $excel.ActiveWorkbook

$oldTimestamp = "01:12:47"
$newTimestamp = "02:38:35"

# --- Sheet 1: LP1912 -------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTimestamp"
$ws1.Range("A3").Value = "Total filas: 3"

# Row 6: 14_ABASTO -> 15_ABASTO (new arrival time/minutes)
$ws1.Cells.Item(6, 1).Value = $newTimestamp
$ws1.Cells.Item(6, 2).Value = "03:01"
$ws1.Cells.Item(6, 3).Value = "15_ABASTO"
$ws1.Cells.Item(6, 4).Value = 23
$ws1.Cells.Item(6, 5).Value = "LP1912"

# Row 7: previously 15_ABASTO, now 14_ABASTO
$ws1.Cells.Item(7, 1).Value = $newTimestamp
$ws1.Cells.Item(7, 2).Value = "03:48"
$ws1.Cells.Item(7, 3).Value = "14_ABASTO"
$ws1.Cells.Item(7, 4).Value = 70
$ws1.Cells.Item(7, 5).Value = "LP1912"

# Row 8: new row added to the schedule
$ws1.Cells.Item(8, 1).Value = $newTimestamp
$ws1.Cells.Item(8, 2).Value = "04:02"
$ws1.Cells.Item(8, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(8, 4).Value = 84
$ws1.Cells.Item(8, 5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTimestamp"

# --- Sheet 3: 6203-6173 ------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTimestamp"
